$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: "Binary Search" / "Special Ques" topic row
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Binary Search"
$ws.Range("C5").Value = "Special Ques"
$ws.Range("D5").Value = "Given a sorted array A that has been rotated in a cycle, find the smallest element of the array in O(log(n)) time. For example,`n[1,2,4,5,7,8] rotated by 3 gives us A = [5,7,8,1,2,4] and the smallest number is 1.`n[1,2,4,5,7,8] rotated by 1 gives us A = [8,1,2,4,5,7] and the smallest number is 1."
$ws.Range("D5").WrapText = $true

$ws.Hyperlinks.Add($ws.Range("E5"), "https://interviewcamp.io/courses/101687/lectures/2636568", "", "", "https://interviewcamp.io/courses/101687/lectures/2636568")

$ws.Rows.Item(5).RowHeight = 85

$ws.Range("E5").Select()
